$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This workbook previously mixed per-row stock data from many different
# tickers (fixed_ticker column) into a file meant to hold only Advanced
# Micro Devices (AMD) records. Replace the open/close/high/low prices and
# shares_outstanding for every data row (2-43) with the correct AMD-only
# values, and normalize fixed_ticker to "AMD" for every row so the extra,
# now-unused ticker strings are dropped from the shared string table.

$amdData = @(
    @(2, 2.660000085830688, 2.259999990463257, 2.940000057220459, 2.25),
    @(3, 2.470000028610229, 1.929999947547913, 2.630000114440918, 1.610000014305115),
    @(4, 1.769999980926514, 2.119999885559082, 2.210000038146973, 1.720000028610229),
    @(5, 2.769999980926514, 2.200000047683716, 2.819999933242798, 1.75),
    @(6, 2.789999961853028, 3.549999952316284, 3.990000009536743, 2.599999904632568),
    @(7, 5.090000152587891, 6.860000133514404, 7.159999847412109, 4.820000171661377),
    @(8, 6.949999809265137, 7.230000019073486, 7.53000020980835, 6.239999771118164),
    @(9, 11.42000007629394, 10.36999988555908, 11.6899995803833, 9.420000076293944),
    @(10, 14.60000038146973, 13.30000019073486, 14.73999977111816, 12.22000026702881),
    @(11, 12.56999969482422, 13.60999965667725, 15.64999961853027, 12.13000011444092),
    @(12, 12.80000019073486, 10.98999977111816, 14.40999984741211, 10.64999961853027),
    @(13, 10.42000007629394, 13.73999977111816, 13.85000038146973, 10.34000015258789),
    @(14, 9.989999771118164, 10.88000011444092, 11.35999965667725, 9.039999961853027),
    @(15, 14.80000019073486, 18.32999992370605, 20.18000030517578, 14.73999961853027),
    @(16, 30.69000053405762, 18.20999908447266, 31.90999984741211, 16.17000007629395),
    @(17, 18.01000022888184, 24.40999984741211, 25.13999938964844, 16.94000053405762),
    @(18, 26.42000007629395, 27.6299991607666, 29.95000076293945, 25.82999992370605),
    @(19, 31.79000091552734, 30.45000076293945, 34.86000061035156, 30.10000038146973),
    @(20, 29.04999923706055, 33.93000030517578, 34.34000015258789, 27.43000030517578),
    @(21, 46.86000061035156, 47, 52.81000137329102, 46.09999847412109),
    @(22, 44.18000030517578, 52.38999938964844, 58.63000106811523, 41.70000076293945),
    @(23, 52.63000106811523, 77.43000030517578, 78.95999908447266, 51.59999847412109),
    @(24, 83.05999755859375, 75.29000091552734, 88.72000122070312, 74.23000335693359),
    @(25, 92.11000061035156, 85.63999938964844, 99.23000335693359, 85.01999664306641),
    @(26, 80.16000366210938, 81.62000274658203, 89.19999694824219, 77.94000244140625),
    @(27, 94.04000091552734, 106.1900024414062, 106.9700012207031, 84.23999786376953),
    @(28, 102.5999984741211, 120.2300033569336, 128.0800018310547, 99.81999969482422),
    @(29, 145.1399993896484, 114.25, 152.4199981689453, 99.34999847412109),
    @(30, 110.4800033569336, 85.51999664306641, 111.4199981689453, 84.01999664306641),
    @(31, 75.19000244140625, 94.47000122070312, 94.80999755859376, 71.59999847412109),
    @(32, 64.45999908447266, 60.06000137329102, 70.29000091552734, 54.56999969482422),
    @(33, 66, 75.15000152587891, 77.08000183105469, 60.04999923706055),
    @(34, 96.6999969482422, 89.37000274658203, 97.26999664306641, 83.76000213623047),
    @(35, 115.1600036621094, 114.4000015258789, 122.120002746582, 108.5500030517578),
    @(36, 102.2099990844727, 98.5, 111.3099975585938, 93.12000274658205),
    @(37, 144.2799987792969, 167.6900024414062, 184.9199981689453, 133.7400054931641),
    @(38, 180.1000061035156, 158.3800048828125, 187.2400054931641, 145.2899932861328),
    @(39, 161.25, 144.4799957275391, 187.2799987792969, 134.0500030517578),
    @(40, 164.4799957275391, 144.0700073242188, 174.0500030517578, 143.3300018310547),
    @(41, 122.2900009155273, 115.9499969482422, 131.7100067138672, 112.8000030517578),
    @(42, 102.1600036621094, 97.34999847412109, 104.1900024414062, 76.48000335693359),
    @(43, 138.7700042724609, 176.3099975585938, 182.5, 133.5)
)

foreach ($row in $amdData) {
    $r = $row[0]
    $ws.Range("D$r").Value = $row[1]
    $ws.Range("E$r").Value = $row[2]
    $ws.Range("F$r").Value = $row[3]
    $ws.Range("G$r").Value = $row[4]
    $ws.Range("H$r").Value = 1622843689
    $ws.Range("I$r").Value = "AMD"
}
